$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URI")

# Row 6: "Change in inventories"
$ws.Range("B6").Value = 109000000.0
$ws.Range("C6").Value = 103000000.0
$ws.Range("D6").Value = 129000000.0
$ws.Range("E6").Value = 29000000.0
$ws.Range("F6").Value = 6000000.0
$ws.Range("G6").Value = -8000000.0

# Row 8: "Change in payables and accrued liability"
$ws.Range("B8").Value = 1302000000.0
$ws.Range("C8").Value = 1807000000.0
$ws.Range("D8").Value = 1795000000.0
$ws.Range("E8").Value = 1300000000.0
$ws.Range("F8").Value = 1026000000.0
$ws.Range("G8").Value = 468000000.0
